$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.135.06"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.557.05"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.41"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.62"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").Value = "3.557.13"
$ws.Range("E7").Value = "  +4.17%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +4.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.92"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "4.154.78"
$ws.Range("E13").Value = "  +3.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000207"
$ws.Range("E14").Value = "  +3.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.08"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "3.556.99"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").Value = "66.207.38"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.115"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.32"
$ws.Range("E19").Value = "  +9.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.41"
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.612"
$ws.Range("E23").Value = "  +6.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.22"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "3.692.67"
$ws.Range("E25").Value = "  +3.98%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  +8.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.08"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.47"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.55"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").Value = "3.546.19"
$ws.Range("E34").Value = "  +3.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.152"
$ws.Range("E35").Value = "  -4.66%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.75"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.89"
$ws.Range("E38").Value = "  +4.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.62"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.58"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0851"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.892"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.93"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.08"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.73"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.65"
$ws.Range("E49").Value = "  +15.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.13"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("E51").Value = "  +2.90%  "
